$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.749.53'
$ws.Range('E2').Value = '  +1.08%  '

$ws.Range('D3').Value = '2.497.08'
$ws.Range('E3').Value = '  +0.83%  '

$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '587.47'
$ws.Range('E5').Value = '  +0.63%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '177.16'
$ws.Range('E6').Value = '  +4.62%  '

$ws.Range('E7').Value = '  -0.03%  '

$ws.Range('E8').Value = '  +0.97%  '

$ws.Range('E9').Value = '  +5.01%  '

$ws.Range('E10').Value = '  +0.46%  '

$ws.Range('E11').Value = '  +3.56%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.94'
$ws.Range('E12').Value = '  +0.73%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '25.78'
$ws.Range('E13').Value = '  +2.10%  '

$ws.Range('D14').Value = '2.917.26'
$ws.Range('E14').Value = '  -0.18%  '

$ws.Range('D15').Value = '67.675.37'
$ws.Range('E15').Value = '  +1.10%  '

$ws.Range('E16').Value = '  +2.42%  '

$ws.Range('D17').Value = '2.506.65'
$ws.Range('E17').Value = '  +0.63%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.08'
$ws.Range('E18').Value = '  +1.11%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.52'
$ws.Range('E19').Value = '  +2.09%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '351.10'
$ws.Range('E20').Value = '  +0.45%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.09'
$ws.Range('E21').Value = '  +2.77%  '

$ws.Range('E22').Value = '  +0.06%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.65'
$ws.Range('E23').Value = '  +3.29%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.28'
$ws.Range('E24').Value = '  +1.91%  '

$ws.Range('E25').Value = '  -1.24%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.17'
$ws.Range('E26').Value = '  +0.65%  '

$ws.Range('E27').Value = '  +1.12%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.998'
$ws.Range('E28').Value = '  -0.09%  '

$ws.Range('D29').Value = '0.0₃0911'
$ws.Range('E29').Value = '  +1.78%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '507.65'
$ws.Range('E30').Value = '  +0.31%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.81'
$ws.Range('E31').Value = '  +2.69%  '

$ws.Range('E32').Value = '  +3.38%  '

$ws.Range('E33').Value = '  +1.08%  '

$ws.Range('E34').Value = '  +0.04%  '

$ws.Range('E35').Value = '  +7.21%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '162.36'
$ws.Range('E36').Value = '  +2.43%  '

$ws.Range('B37').Value = 'WhiteBITCoin'
$ws.Range('C37').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.67'
$ws.Range('E37').Value = '  +0.04%  '

$ws.Range('B38').Value = 'EthereumClassic'
$ws.Range('C38').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.42'
$ws.Range('E38').Value = '  +1.42%  '

$ws.Range('E39').Value = '  +1.57%  '

$ws.Range('E40').Value = '  +0.01%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.74'
$ws.Range('E41').Value = '  +4.27%  '

$ws.Range('E42').Value = '  +1.58%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.87'
$ws.Range('E43').Value = '  +1.78%  '

$ws.Range('E44').Value = '  +3.42%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '145.57'
$ws.Range('E45').Value = '  +3.18%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.51'
$ws.Range('E46').Value = '  +2.57%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.517'
$ws.Range('E47').Value = '  +1.27%  '

$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0744'
$ws.Range('E48').Value = '  +2.57%  '

$ws.Range('B49').Value = 'Optimism'
$ws.Range('C49').Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.59'
$ws.Range('E49').Value = '  +2.21%  '

$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.587'
$ws.Range('E50').Value = '  +1.13%  '

$ws.Range('B51').Value = 'BitgetToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.17'
$ws.Range('E51').Value = '  +0.04%  '
